# Updated via Streamlit Approval System
# Rows 82-86 get new beneficiary-approval data (each row "rotates" into the
# next request), and a brand-new row 87 is appended for the latest request.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 82 ---
$ws.Range("A82").Value = "WGE 203"
$ws.Range("G82").Value = "DCR"
$ws.Range("K82").Value = "MADUMITHA"
$ws.Range("L82").Value = "94b2efed-2699-4f62-92a0-8812abcc9e78"
$ws.Range("M82").Value = "ACC-10629442465"
$ws.Range("N82").Value = "SBIN0002016"
$ws.Range("V82").Value = 15000
$ws.Range("X82").Value = "room rent  jan 26 RPA_ID : 3627b53840"
$ws.Range("Y82").Value = "kolkata"

# --- Row 83 ---
$ws.Range("A83").Value = "WGP011"
$ws.Range("G83").Value = "NEFT"
$ws.Range("K83").Value = "SHREE BALAJI ELECTRICAL"
$ws.Range("L83").Value = "d892dbf3-8741-44f1-ba7f-5a584bc7f350"
$ws.Range("M83").Value = "ACC-125006695576"
$ws.Range("N83").Value = "CNRB0017203"
$ws.Range("V83").Value = 105987
$ws.Range("X83").Value = "Being electric consumables purchased RPA_ID : abb7bb472c"
$ws.Range("Y83").Value = "ONGC Electrical"
$ws.Range("Z83").Value = "SITE EXPENSE"
$ws.Range("AA83").Value = "midhuncraju12@gmail.com"

# --- Row 84 ---
$ws.Range("A84").Value = "WGG 02"
$ws.Range("K84").Value = ""
$ws.Range("L84").Value = "788a71f1-06f3-4161-8e96-7dbdaaa092ca"
$ws.Range("M84").Value = ""
$ws.Range("N84").Value = ""
$ws.Range("V84").Value = 1470
$ws.Range("X84").Value = "Being IOCL Willington switch and core cutting charges RPA_ID : 183a5be1f1"
$ws.Range("Y84").Value = "IOCL Willington"

# --- Row 85 ---
$ws.Range("A85").Value = "WGP008"
$ws.Range("G85").Value = "DCR"
$ws.Range("K85").Value = "GAYATHRI ELECTRICALS"
$ws.Range("L85").Value = "77e4cafa-fced-4e92-af02-1d695f6c561d"
$ws.Range("M85").Value = "ACC-39177475703"
$ws.Range("N85").Value = "SBIN0000512"
$ws.Range("V85").Value = 580870
$ws.Range("X85").Value = "Being material purchase RPA_ID : 32aed1f5ef"
$ws.Range("Y85").Value = "ONGC Electrical"
$ws.Range("Z85").Value = "SITE EXPENSES"

# --- Row 86 ---
$ws.Range("A86").Value = "WGE 234"
$ws.Range("G86").Value = "NEFT"
$ws.Range("K86").Value = "MANU"
$ws.Range("L86").Value = "a82b0afd-403a-4722-9101-1bfb56852a57"
$ws.Range("M86").Value = "ACC-111001506458"
$ws.Range("N86").Value = "ICIC0001110"
$ws.Range("V86").Value = 2000
$ws.Range("X86").Value = "Purchase of stamp papers RPA_ID : d3b8adacc7"
$ws.Range("Y86").Value = "RO site Thiruvaniyoor"
$ws.Range("Z86").Value = "FOR FORM OF CONTRACT"

# --- Row 87 (new) ---
$ws.Range("A87").Value = "WGE 234"
$ws.Range("B87").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("C87").Value = "13-02-2026"
$ws.Range("D87").Value = 286962
$ws.Range("E87").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("F87").Value = 34413429360
$ws.Range("G87").Value = "NEFT"
$ws.Range("H87").Value = "SBIN0003229"
$ws.Range("I87").Value = "AAAFW8862C"
$ws.Range("J87").Value = "32AAAFW8862C1Z9"
$ws.Range("K87").Value = "MANU"
$ws.Range("L87").Value = "8e160b61-efa7-4af3-b7ed-b023930becd4"
$ws.Range("M87").Value = "ACC-111001506458"
$ws.Range("N87").Value = "ICIC0001110"
$ws.Range("U87").Value = "pending"
$ws.Range("V87").Value = 100
$ws.Range("X87").Value = "Petrol expense RPA_ID : 936b48b8c6"
$ws.Range("Y87").Value = "RO site Thiruvaniyoor"
$ws.Range("Z87").Value = "FOR FORM OF CONTRACT"
$ws.Range("AA87").Value = "midhuncraju12@gmail.com"
$ws.Range("AB87").Value = "ESTIMATION NOT MATCHED"
$ws.Range("AC87").Value = 0
$ws.Range("AD87").Value = 0
$ws.Range("AE87").Value = 0
